$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14
$ws.Cells.Item($row, 1).Value = "'2026-01-20"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "'2026-01-20 00:59:45"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Value = 2.15
$ws.Cells.Item($row, 4).Value = 115
$ws.Cells.Item($row, 5).Value = 0.85
$ws.Cells.Item($row, 6).Value = 1102
$ws.Cells.Item($row, 7).Value = 1.3
$ws.Cells.Item($row, 8).Value = 1217
$ws.Cells.Item($row, 9).Value = 2.75
$ws.Cells.Item($row, 10).Value = 2
$ws.Cells.Item($row, 11).Value = 1.22
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 1.53
$ws.Cells.Item($row, 14).Value = 3
